$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D6").Value2 = "[name=`"위스퍼레인`"]  이 이상의 치료는 의료부 전문 의사의 도움을 받아야 해요.`n"
$ws.Range("D24").Value2 = "[name=`"폴리닉`"]  당신은 정식 의료부 닥터는 아니지만, 저로서는 당신에게 충고해 둘 책임이 있습니다.`n"
$ws.Range("C25").Value2 = "[name=`"Folinic`"]  According to the 'Regulations for Management and Treatment of Rhodes Island Wounded Personnel,' all wounded Operators must undergo registration at the Rhodes Island Medical Department.`n"
$ws.Range("D26").Value2 = "[name=`"폴리닉`"]  당신처럼 개인적으로 오퍼레이터 상처를 치료해 주는 것에 대해 작은 상처나 작은 병은 괜찮을지 몰라도, 만약 팬텀 오퍼레이터가 당신의 불완전한 치료로 합병증을 얻거나 상처가 악화하면 의료부는 아주 골치가 아플 겁니다.`n"
$ws.Range("D46").Value2 = "[name=`"폴리닉`"]  긴장하지 마세요. 전 당신을 데리고 로비를 지나 의료부 쪽으로 갈 생각은 없으니까요.`n"
$ws.Range("C61").Value2 = "[name=`"Whisperain`"]  I believe it's because you told him 'mouth shut' just prior, Folinic.`n"
$ws.Range("C99").Value2 = "[name=`"Folinic`"]  'It's lucky we had him supporting us. We wouldn't have come back alive otherwise.' I've already heard from plenty of others to that effect.`n"
$ws.Range("D114").Value2 = "[name=`"팬텀`"]  나중에 블루포이즌에게 직접 감사 인사를 전하도록 하지`n"
$ws.Range("D132").Value2 = "[name=`"폴리닉`"]  물론 가장 중요한 것은 정기적으로 의료부를 찾아와 치료를 받는 거예요.`n"
$ws.Range("D134").Value2 = "[name=`"폴리닉`"]  의료부에 오기 싫다면, 제가 다른 사람과의 접촉을 줄일 수 있는 순회 진료 쪽으로 알아봐 드릴게요.`n"
$ws.Range("C174").Value2 = "[name=`"Phantom`"]  'Wanted: Multiple staff members for a Multimedia Escape Experience event—'`n"
$ws.Range("C176").Value2 = "[name=`"Phantom`"]  'Staff required to take charge in planning escape rooms in whole, for an urgent and pressing adventurous atmosphere.'`n"
$ws.Range("C177").Value2 = "[name=`"Phantom`"]  'Staff required to plan deciphering and mechanisms within escape rooms.'`n"
$ws.Range("C178").Value2 = "[name=`"Phantom`"]  'Staff required to—'`n"
$ws.Range("C179").Value2 = "[name=`"Phantom`"]  '—act out roles.'`n"
$ws.Range("C181").Value2 = "[name=`"Phantom`"]  'Design (until opening) to last one month. Salaried for the duration, with three meals a day included.'`n"
$ws.Range("C182").Value2 = "[name=`"Phantom`"]  'Interested operators and staff aboard, please find details at the Human Resources Department.'`n"
$ws.Range("C189").Value2 = "[name=`"Whisperain`"]  'It scared me to death!'`n"

# Excel auto-expands row height to fit new content on edit; restore the original
# auto (non-custom) row height behavior since the source rows had no custom height.
$ws.Rows.Item(6).EntireRow.AutoFit()
$ws.Rows.Item(24).EntireRow.AutoFit()
$ws.Rows.Item(25).EntireRow.AutoFit()
$ws.Rows.Item(26).EntireRow.AutoFit()
$ws.Rows.Item(46).EntireRow.AutoFit()
$ws.Rows.Item(61).EntireRow.AutoFit()
$ws.Rows.Item(99).EntireRow.AutoFit()
$ws.Rows.Item(114).EntireRow.AutoFit()
$ws.Rows.Item(132).EntireRow.AutoFit()
$ws.Rows.Item(134).EntireRow.AutoFit()
$ws.Rows.Item(174).EntireRow.AutoFit()
$ws.Rows.Item(176).EntireRow.AutoFit()
$ws.Rows.Item(177).EntireRow.AutoFit()
$ws.Rows.Item(178).EntireRow.AutoFit()
$ws.Rows.Item(179).EntireRow.AutoFit()
$ws.Rows.Item(181).EntireRow.AutoFit()
$ws.Rows.Item(182).EntireRow.AutoFit()
$ws.Rows.Item(189).EntireRow.AutoFit()
